$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day 11 (row 15): new Part 1 time, Part 2 becomes a plain recorded value,
# and Total becomes the SUM formula (matches the shared formula column E uses).
$ws.Range("B15").Value = 0.0533150000264868
$ws.Range("C15").Value = 0.000448399921879172
$ws.Range("E15").Formula = "=SUM(B15:C15)"

# Day 12 (row 16): now has its own runtimes; Part2 uses the +E-B formula style
# (matching the one previously used in row 15), Total becomes a plain value.
$ws.Range("B16").Value = 0.0115254999836906
$ws.Range("C16").Formula = "=+E16-B16"
$ws.Range("E16").Value = 0.187391299987211

# Day 13 (row 17): newly filled in runtimes, Total keeps its formula.
$ws.Range("B17").Value = 0.00746969995088875
$ws.Range("C17").Value = 0.000111200031824409

$ws.Range("E17").Select()
